$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Rewrite the changed sentence inside the existing intro paragraph.
#    "...esa energia mecanica generada por las personas transformandola
#     en energia renovable y limpia aprovechando lo maximo posible el
#     transito, esto con el objetivo de poder mantener..."
#    becomes
#    "...ese movimiento transformandolo en energia aprovechando al
#     maximo el transito para mantener..."
# ---------------------------------------------------------------------
$old1 = "esa energía mecánica generada por las personas transformándola en energía renovable y limpia aprovechando lo máximo posible el tránsito, esto con el objetivo de poder mantener"
$new1 = "ese movimiento transformándolo en energía aprovechando al máximo el tránsito para mantener"
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Append the new sections after the (now single-run) intro paragraph.
#    Each new paragraph is created as a *plain* paragraph mark first
#    (inheriting the preceding plain body style), and only afterwards
#    is it given heading formatting (centered, bold, size 36) when it
#    is meant to be a heading - this keeps the following paragraph from
#    inheriting the bold/centered run formatting.
# ---------------------------------------------------------------------
$introPara = $d.Paragraphs(2)
$r = $introPara.Range
$r.Collapse(0)

# --- paragraph 3: heading "¿Por qué se utilizaría?" ---
$r.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$r = $p3.Range
$r.Collapse(0)

# --- paragraph 4: body text ---
$r.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$r = $p4.Range
$r.Collapse(0)

# --- paragraph 5: heading "Explicación piezoelectricidad" ---
$r.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$r = $p5.Range
$r.Collapse(0)

# --- paragraph 6: empty paragraph holding the _GoBack bookmark ---
$r.InsertParagraphAfter()
$p6 = $d.Paragraphs(6)
$p6.Range.Bookmarks.Add("_GoBack") | Out-Null
$r = $p6.Range
$r.Collapse(0)

# --- paragraph 7: heading "¿Cómo funciona la baldosa?" ---
$r.InsertParagraphAfter()
$p7 = $d.Paragraphs(7)
$r = $p7.Range
$r.Collapse(0)

# --- paragraph 8: body text ---
$r.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)

# -----------------------------------------------------------------
# Now fill in text + formatting, heading paragraphs last so their
# bold/centered formatting never leaks onto a still-to-be-created
# sibling paragraph.
# -----------------------------------------------------------------
$p4.Range.Text = "Los carteles de salida de emergencias deben estar colocados en cada salida de forma obligatoria colocadas a simple vistas para ser reconocidas rápidamente para reducir lo máximo posible cualquier riesgo, por lo que en lugares sin luz natural como las estaciones de subte deben tener un sistema de iluminación autónomo que las mantenga encendidas incluso en un corte de luz. Pleper se crea como una forma de mantener cargados los carteles de salida para evitar problemas como la descarga con el tiempo de estos mismos, ayudando a reducir el riesgo y funcionando como una alternativa limpia y renovable de aprovechar al máximo la energía generada por las personas."

$p8.Range.Text = "La baldosa consiste en una plancha de madera apoyada sobre cuatro resortes que funcionan como pilares y permiten bajar la baldosa hasta deformar los piezoeléctricos y luego regresarlo a su posición original. Bajo la plancha de madera se encuentran los piezoeléctricos divididos en 8 grupos de 4 piezoeléctricos colocados en paralelos y rectificados en la placa"

$p3.Range.Text = "¿Por qué se utilizaría?"
$p3.Alignment = 1
$p3.Range.Font.Bold = 1
$p3.Range.Font.Size = 18

$p5.Range.Text = "Explicación piezoelectricidad"
$p5.Alignment = 1
$p5.Range.Font.Bold = 1
$p5.Range.Font.Size = 18

$p7.Range.Text = "¿Cómo funciona la baldosa?"
$p7.Alignment = 1
$p7.Range.Font.Bold = 1
$p7.Range.Font.Size = 18

Write-Output "Done. Paragraph count:"
Write-Output $d.Paragraphs.Count
